$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value2 = 762.1429000000001
$ws.Range("I55").Value2 = 1003.25
$ws.Range("J55").Value2 = 440.66666
$ws.Range("K55").Value2 = 1003.25
$ws.Range("L55").Value2 = 440.66666
$ws.Range("M55").Value2 = -789.25
$ws.Range("N55").Value2 = -868.66666

$ws.Range("H68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("L68").Value2 = 0
$ws.Range("N68").ClearContents()

$ws.Range("H69").Value2 = 7671.143
$ws.Range("I69").Value2 = 6739.6
$ws.Range("K69").Value2 = 20218.8
$ws.Range("M69").Value2 = -19344.8

$ws.Range("H71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("L71").Value2 = 0
$ws.Range("N71").ClearContents()

$ws.Range("H72").Value2 = 7671.143
$ws.Range("I72").Value2 = 6739.6
$ws.Range("K72").Value2 = 60656.4
$ws.Range("M72").Value2 = -56288.4

$ws.Range("H80").Value2 = 337.6
$ws.Range("I80").Value2 = 359.5
$ws.Range("K80").Value2 = 1078.5
$ws.Range("M80").Value2 = -80.5

$ws.Range("H83").Value2 = 337.6
$ws.Range("I83").Value2 = 359.5
$ws.Range("K83").Value2 = 3235.5
$ws.Range("M83").Value2 = 1756.5

$ws.Range("H138").Value2 = 3277.7847
$ws.Range("I138").Value2 = 3235.875
$ws.Range("J138").Value2 = 3302.3171
$ws.Range("K138").Value2 = 9707.625
$ws.Range("L138").Value2 = 9906.951300000001
$ws.Range("M138").Value2 = -4567.625
$ws.Range("N138").Value2 = -20186.9513

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value2 = 1250
$ws.Range("I19").Value2 = 1250
$ws.Range("K19").Value2 = 1250
$ws.Range("M19").Value2 = -1021

$ws.Range("H61").Value2 = 2648.75
$ws.Range("I61").Value2 = 2827.2856
$ws.Range("K61").Value2 = 2827.2856
$ws.Range("M61").Value2 = -2615.2856

$ws.Range("H136").Value2 = 2648.75
$ws.Range("I136").Value2 = 2827.2856
$ws.Range("K136").Value2 = 8481.856800000001
$ws.Range("M136").Value2 = -5931.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value2 = 44895.332
$ws.Range("J81").Value2 = 44895.332
$ws.Range("L81").Value2 = 44895.332
$ws.Range("N81").Value2 = -47017.332

$ws.Range("H84").Value2 = 44895.332
$ws.Range("J84").Value2 = 44895.332
$ws.Range("L84").Value2 = 134685.996
$ws.Range("N84").Value2 = -145293.996

$ws.Range("H134").Value2 = 3996.25
$ws.Range("I134").Value2 = 3996.25
$ws.Range("K134").Value2 = 11988.75
$ws.Range("M134").Value2 = -9453.75

$ws.Range("H138").Value2 = 550000.5
$ws.Range("J138").Value2 = 550000.5
$ws.Range("L138").Value2 = 550000.5
$ws.Range("N138").Value2 = -560280.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 840.4167
$ws.Range("I16").Value2 = 819.8570999999999
$ws.Range("J16").Value2 = 869.2
$ws.Range("K16").Value2 = 819.8570999999999
$ws.Range("L16").Value2 = 869.2
$ws.Range("M16").Value2 = -532.8570999999999
$ws.Range("N16").Value2 = -1443.2

$ws.Range("H31").Value2 = 2869.6785
$ws.Range("I31").Value2 = 2519.9
$ws.Range("K31").Value2 = 2519.9
$ws.Range("M31").Value2 = -2224.9

$ws.Range("H34").Value2 = 2869.6785
$ws.Range("I34").Value2 = 2519.9
$ws.Range("K34").Value2 = 2519.9
$ws.Range("M34").Value2 = -2317.9

$ws.Range("H68").Value2 = 85799.164
$ws.Range("J68").Value2 = 85799.164
$ws.Range("L68").Value2 = 85799.164
$ws.Range("N68").Value2 = -87297.164

$ws.Range("H71").Value2 = 85799.164
$ws.Range("J71").Value2 = 85799.164
$ws.Range("L71").Value2 = 257397.492
$ws.Range("N71").Value2 = -264885.492

$ws.Range("H99").Value2 = 0
$ws.Range("I99").Value2 = 0
$ws.Range("J99").Value2 = 0
$ws.Range("K99").Value2 = 0
$ws.Range("L99").Value2 = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H113").Value2 = 840.4167
$ws.Range("I113").Value2 = 819.8570999999999
$ws.Range("J113").Value2 = 869.2
$ws.Range("K113").Value2 = 819.8570999999999
$ws.Range("L113").Value2 = 869.2
$ws.Range("M113").Value2 = 1350.1429
$ws.Range("N113").Value2 = -5209.2

$ws.Range("H126").Value2 = 0
$ws.Range("I126").Value2 = 0
$ws.Range("J126").Value2 = 0
$ws.Range("K126").Value2 = 0
$ws.Range("L126").Value2 = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value2 = 288
$ws.Range("I14").Value2 = 288
$ws.Range("K14").Value2 = 864
$ws.Range("M14").Value2 = -691

$ws.Range("H25").Value2 = 313.8
$ws.Range("J25").Value2 = 0
$ws.Range("L25").Value2 = 0
$ws.Range("N25").ClearContents()

$ws.Range("H30").Value2 = 313.8
$ws.Range("J30").Value2 = 0
$ws.Range("L30").Value2 = 0
$ws.Range("N30").ClearContents()

$ws.Range("H80").Value2 = 50230
$ws.Range("I80").Value2 = 58667.332
$ws.Range("J80").Value2 = 47417.555
$ws.Range("K80").Value2 = 176001.996
$ws.Range("L80").Value2 = 142252.665
$ws.Range("M80").Value2 = -175065.996
$ws.Range("N80").Value2 = -144124.665

$ws.Range("H83").Value2 = 50230
$ws.Range("I83").Value2 = 58667.332
$ws.Range("J83").Value2 = 47417.555
$ws.Range("K83").Value2 = 528005.988
$ws.Range("L83").Value2 = 426757.995
$ws.Range("M83").Value2 = -523325.988
$ws.Range("N83").Value2 = -436117.995

$ws.Range("H97").Value2 = 3875
$ws.Range("J97").Value2 = 1500
$ws.Range("L97").Value2 = 4500
$ws.Range("N97").Value2 = -5492

$ws.Range("H107").Value2 = 1020.6667
$ws.Range("J107").Value2 = 1104.9286
$ws.Range("L107").Value2 = 3314.7858
$ws.Range("N107").Value2 = -7154.7858

$ws.Range("H109").Value2 = 3242
$ws.Range("I109").Value2 = 2365.2
$ws.Range("K109").Value2 = 7095.599999999999
$ws.Range("M109").Value2 = -6055.599999999999

$ws.Range("H132").Value2 = 1525.6666
$ws.Range("J132").Value2 = 1488.5
$ws.Range("L132").Value2 = 13396.5
$ws.Range("N132").Value2 = -18456.5

$ws.Range("H139").Value2 = 14026
$ws.Range("I139").Value2 = 10080.333
$ws.Range("K139").Value2 = 30240.999
$ws.Range("M139").Value2 = -25100.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value2 = 49899.5
$ws.Range("J26").Value2 = 49899.5
$ws.Range("L26").Value2 = 49899.5
$ws.Range("N26").Value2 = -50459.5

$ws.Range("H50").Value2 = 49899.5
$ws.Range("J50").Value2 = 49899.5
$ws.Range("L50").Value2 = 49899.5
$ws.Range("N50").Value2 = -50895.5

$ws.Range("H70").Value2 = 9995.223
$ws.Range("I70").Value2 = 9094.5
$ws.Range("K70").Value2 = 9094.5
$ws.Range("M70").Value2 = -8824.5

$ws.Range("H73").Value2 = 9995.223
$ws.Range("I73").Value2 = 9094.5
$ws.Range("K73").Value2 = 9094.5
$ws.Range("M73").Value2 = -8158.5

$ws.Range("H104").Value2 = 88888
$ws.Range("J104").Value2 = 88888
$ws.Range("L104").Value2 = 88888
$ws.Range("N104").Value2 = -95876

$ws.Range("H132").Value2 = 1759.5
$ws.Range("I132").Value2 = 1759.5
$ws.Range("K132").Value2 = 5278.5
$ws.Range("M132").Value2 = -2748.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1577.1
$ws.Range("I22").Value2 = 881.3333
$ws.Range("K22").Value2 = 881.3333
$ws.Range("M22").Value2 = -586.3333

$ws.Range("H27").Value2 = 1577.1
$ws.Range("I27").Value2 = 881.3333
$ws.Range("K27").Value2 = 881.3333
$ws.Range("M27").Value2 = -774.3333

$ws.Range("H55").Value2 = 631.6667
$ws.Range("J55").Value2 = 617
$ws.Range("L55").Value2 = 617
$ws.Range("N55").Value2 = -963

$ws.Range("H61").Value2 = 1993.3334
$ws.Range("I61").Value2 = 1992.5
$ws.Range("K61").Value2 = 1992.5
$ws.Range("M61").Value2 = -1790.5

$ws.Range("H68").Value2 = 2856.4666
$ws.Range("I68").Value2 = 2567.9092
$ws.Range("J68").Value2 = 3650
$ws.Range("K68").Value2 = 2567.9092
$ws.Range("L68").Value2 = 3650
$ws.Range("M68").Value2 = -1818.9092
$ws.Range("N68").Value2 = -5148

$ws.Range("H71").Value2 = 2856.4666
$ws.Range("I71").Value2 = 2567.9092
$ws.Range("J71").Value2 = 3650
$ws.Range("K71").Value2 = 12839.546
$ws.Range("L71").Value2 = 18250
$ws.Range("M71").Value2 = -9095.546
$ws.Range("N71").Value2 = -25738

$ws.Range("H113").Value2 = 1993.3334
$ws.Range("I113").Value2 = 1992.5
$ws.Range("K113").Value2 = 1992.5
$ws.Range("M113").Value2 = 177.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("L63").Value2 = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value2 = 0
$ws.Range("J66").Value2 = 0
$ws.Range("L66").Value2 = 0
$ws.Range("N66").ClearContents()

$ws.Range("H104").Value2 = 15499.75
$ws.Range("J104").Value2 = 15499.75
$ws.Range("L104").Value2 = 15499.75
$ws.Range("N104").Value2 = -22487.75

$ws.Range("H132").Value2 = 2187.56
$ws.Range("I132").Value2 = 2186.5652
$ws.Range("K132").Value2 = 6559.6956
$ws.Range("M132").Value2 = -4029.6956

$ws.Range("H136").Value2 = 1457.3182
$ws.Range("I136").Value2 = 1207.1333
$ws.Range("K136").Value2 = 3621.3999
$ws.Range("M136").Value2 = -1071.3999
